$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have extra "damage_subtype" and "number" columns (D and
# E) that were always populated with the constant value 1. Remove them,
# which shifts the former "type" (F) and "Material" (G) columns left into
# D and E.
$ws.Range("D1:E1").EntireColumn.Delete()

# Row 7 (pipe 225 now, previously was a duplicate-looking row of pipe 204)
# was corrected: pipe id, damage location and damage type values changed.
$ws.Range("B7").Value = 225
$ws.Range("C7").Value = 0.5
$ws.Range("D7").Value = "leak"

# Leave the selection where the author left it when saving the file.
$ws.Range("C12").Select() | Out-Null
